# Insert a new row before row 4 on the first worksheet (shifts rows 4:11 -> 5:12),
# then populate the newly inserted row 4 with the new
# "climate_change_factor_gnrl_hydropower_availability" variable entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a blank row at position 4, pushing existing rows 4-11 down to 5-12.
$ws.Rows.Item(4).Insert()

# Column A: subsector label
$ws.Cells.Item(4, 1).Value = "General"

# Column B: variable name
$ws.Cells.Item(4, 2).Value = "climate_change_factor_gnrl_hydropower_availability"

# Columns C-G stay blank (normalize_group, trajgroup_no_vary_q, uniform_scaling_q,
# variable_trajectory_group, variable_trajectory_group_trajectory_type)

# Column H: max_35
$ws.Cells.Item(4, 8).Value = 1

# Column I: min_35
$ws.Cells.Item(4, 9).Value = 0.5

# Columns J:AS (time series values 0..35) all set to 1
for ($col = 10; $col -le 45; $col++) {
    $ws.Cells.Item(4, $col).Value = 1
}
